$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F1").Value = "SuccessMsgForRecordCreation"
$ws.Range("F5").Value = "SuccessMsgForRecordCreation"

$ws.Range("F7").Select()
